$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 26, shifting existing rows 26-37 down to 27-38.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly data point.
$ws.Cells.Item(26, 1).Value = 10
$ws.Cells.Item(26, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(26, 3).Value = "La Araucanía"
$ws.Cells.Item(26, 4).Value = (Get-Date -Year 2021 -Month 9 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(26, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(26, 5).Value = 9
$ws.Cells.Item(26, 6).Value = 300000001
$ws.Cells.Item(26, 7).Value = "Rabanito"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 50
$ws.Cells.Item(26, 11).Value = 6000
$ws.Cells.Item(26, 12).Value = 7000
$ws.Cells.Item(26, 13).Value = 6400
$ws.Cells.Item(26, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(26, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(26, 16).Value = 533
$ws.Cells.Item(26, 17).Value = 12
$ws.Cells.Item(26, 18).Value = "Hortaliza"
